$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)

# Title shape: remove the "Спасибо за внимание" run, leaving the trailing empty run untouched.
$titleShape = $s.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Runs(1).Text = ""

# Subtitle textbox: replace "Ваши вопросы ?" with "Спасибо за внимание" (typed as two runs).
$subShape = $s.Shapes.Item(2)
$subRange = $subShape.TextFrame.TextRange
$subRange.Text = "Спасибо за "
$null = $subRange.InsertAfter("внимание")
